# Auto-generated edit script: Add data for 2024-10-29
# Updates column K (year 2024 totals) on the Citywide Totals, By Neighborhood,
# and individual neighborhood sheets to reflect one additional day of data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 11).Value = 6730
$ws.Cells.Item(3, 11).Value = 6936
$ws.Cells.Item(4, 11).Value = 1436
$ws.Cells.Item(6, 11).Value = 7614
$ws.Cells.Item(7, 11).Value = 23216

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(2, 11).Value = 75
$ws.Cells.Item(7, 11).Value = 298

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 11).Value = 422
$ws.Cells.Item(3, 11).Value = 464
$ws.Cells.Item(6, 11).Value = 504
$ws.Cells.Item(7, 11).Value = 1523

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(3, 11).Value = 178
$ws.Cells.Item(7, 11).Value = 501

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 11).Value = 356
$ws.Cells.Item(6, 11).Value = 316
$ws.Cells.Item(7, 11).Value = 1002

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 11).Value = 132
$ws.Cells.Item(7, 11).Value = 383

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 11).Value = 259
$ws.Cells.Item(6, 11).Value = 235
$ws.Cells.Item(7, 11).Value = 787

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 11).Value = 179
$ws.Cells.Item(6, 11).Value = 198
$ws.Cells.Item(7, 11).Value = 542

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(6, 11).Value = 99
$ws.Cells.Item(7, 11).Value = 393

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 11).Value = 203
$ws.Cells.Item(6, 11).Value = 164
$ws.Cells.Item(7, 11).Value = 704
$ws.Cells.Item(8, 11).Value = 1523
$ws.Cells.Item(9, 11).Value = 105
$ws.Cells.Item(11, 11).Value = 430
$ws.Cells.Item(15, 11).Value = 245
$ws.Cells.Item(19, 11).Value = 680
$ws.Cells.Item(20, 11).Value = 562
$ws.Cells.Item(21, 11).Value = 75
$ws.Cells.Item(22, 11).Value = 73
$ws.Cells.Item(27, 11).Value = 218
$ws.Cells.Item(29, 11).Value = 1262
$ws.Cells.Item(31, 11).Value = 255
$ws.Cells.Item(33, 11).Value = 1002
$ws.Cells.Item(37, 11).Value = 787
$ws.Cells.Item(48, 11).Value = 295
$ws.Cells.Item(49, 11).Value = 126
$ws.Cells.Item(52, 11).Value = 616
$ws.Cells.Item(53, 11).Value = 298
$ws.Cells.Item(54, 11).Value = 455
$ws.Cells.Item(55, 11).Value = 249
$ws.Cells.Item(57, 11).Value = 87
$ws.Cells.Item(60, 11).Value = 134
$ws.Cells.Item(63, 11).Value = 62
$ws.Cells.Item(64, 11).Value = 145
$ws.Cells.Item(65, 11).Value = 542
$ws.Cells.Item(70, 11).Value = 40
$ws.Cells.Item(73, 11).Value = 209
$ws.Cells.Item(75, 11).Value = 71
$ws.Cells.Item(77, 11).Value = 159
$ws.Cells.Item(78, 11).Value = 266
$ws.Cells.Item(83, 11).Value = 501
$ws.Cells.Item(85, 11).Value = 1069
$ws.Cells.Item(86, 11).Value = 144
$ws.Cells.Item(88, 11).Value = 250
$ws.Cells.Item(89, 11).Value = 348
$ws.Cells.Item(92, 11).Value = 86
$ws.Cells.Item(94, 11).Value = 311
$ws.Cells.Item(95, 11).Value = 383
$ws.Cells.Item(97, 11).Value = 182
$ws.Cells.Item(99, 11).Value = 393
$ws.Cells.Item(101, 11).Value = 23216

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 11).Value = 84
$ws.Cells.Item(7, 11).Value = 255

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(2, 11).Value = 28
$ws.Cells.Item(7, 11).Value = 126

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(6, 11).Value = 245
$ws.Cells.Item(7, 11).Value = 455

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 11).Value = 356
$ws.Cells.Item(3, 11).Value = 447
$ws.Cells.Item(6, 11).Value = 368
$ws.Cells.Item(7, 11).Value = 1262

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 11).Value = 45
$ws.Cells.Item(7, 11).Value = 295

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(3, 11).Value = 204
$ws.Cells.Item(6, 11).Value = 226
$ws.Cells.Item(7, 11).Value = 680

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(2, 11).Value = 62
$ws.Cells.Item(7, 11).Value = 164

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(3, 11).Value = 66
$ws.Cells.Item(6, 11).Value = 91
$ws.Cells.Item(7, 11).Value = 266

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(3, 11).Value = 73
$ws.Cells.Item(7, 11).Value = 249

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Cells.Item(3, 11).Value = 19
$ws.Cells.Item(7, 11).Value = 75

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(3, 11).Value = 40
$ws.Cells.Item(7, 11).Value = 145

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 11).Value = 195
$ws.Cells.Item(3, 11).Value = 181
$ws.Cells.Item(7, 11).Value = 562

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 11).Value = 233
$ws.Cells.Item(6, 11).Value = 192
$ws.Cells.Item(7, 11).Value = 704

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(2, 11).Value = 79
$ws.Cells.Item(6, 11).Value = 141
$ws.Cells.Item(7, 11).Value = 311

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(6, 11).Value = 72
$ws.Cells.Item(7, 11).Value = 245

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 11).Value = 151
$ws.Cells.Item(7, 11).Value = 430

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(4, 11).Value = 6
$ws.Cells.Item(7, 11).Value = 105

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 11).Value = 72
$ws.Cells.Item(7, 11).Value = 209

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(3, 11).Value = 56
$ws.Cells.Item(7, 11).Value = 203

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(6, 11).Value = 97
$ws.Cells.Item(7, 11).Value = 182

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Cells.Item(2, 11).Value = 26
$ws.Cells.Item(7, 11).Value = 86

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Cells.Item(2, 11).Value = 19
$ws.Cells.Item(7, 11).Value = 40

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(3, 11).Value = 77
$ws.Cells.Item(7, 11).Value = 250

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(2, 11).Value = 97
$ws.Cells.Item(6, 11).Value = 103
$ws.Cells.Item(7, 11).Value = 348

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(3, 11).Value = 51
$ws.Cells.Item(7, 11).Value = 218

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(2, 11).Value = 25
$ws.Cells.Item(4, 11).Value = 62
$ws.Cells.Item(7, 11).Value = 144

$ws = $wb.Worksheets.Item('Pullman')
$ws.Cells.Item(6, 11).Value = 13
$ws.Cells.Item(7, 11).Value = 71

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(6, 11).Value = 38
$ws.Cells.Item(7, 11).Value = 87

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(3, 11).Value = 41
$ws.Cells.Item(7, 11).Value = 134

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 11).Value = 352
$ws.Cells.Item(3, 11).Value = 372
$ws.Cells.Item(6, 11).Value = 260
$ws.Cells.Item(7, 11).Value = 1069

$ws = $wb.Worksheets.Item('Clearing')
$ws.Cells.Item(4, 11).Value = 4
$ws.Cells.Item(7, 11).Value = 73

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Cells.Item(2, 11).Value = 67
$ws.Cells.Item(7, 11).Value = 159

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 11).Value = 165
$ws.Cells.Item(6, 11).Value = 224
$ws.Cells.Item(7, 11).Value = 616
